$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for the 3x3 sending/target cluster matrix (ECs, FAPs, sCs)
$data = @(
    @("ECs",  "Sema6a", "Plxna4", "ECs",  2, 0.6666666666666666, 42.665376,          127.996128,          0.5598448706913429, 0.5598448706913429, 2, 0.6666666666666666, 2.431926, 7.295778,           0.8146397085966534, 0.8146397085966534, 103.759037194176,  933.831334747584,  0.4560718623193267,  0.4560718623193267),
    @("ECs",  "Sema6a", "Plxna4", "FAPs", 2, 0.6666666666666666, 42.665376,          127.996128,          0.5598448706913429, 0.5598448706913429, 3, 1,                  0.220323, 0.6609689999999999, 0.07380317678956531, 0.07380317678956533, 9.400163636447999, 84.60147272803199, 0.04131832996636452, 0.04131832996636452),
    @("ECs",  "Sema6a", "Plxna4", "sCs",  2, 0.6666666666666666, 42.665376,          127.996128,          0.5598448706913429, 0.5598448706913429, 3, 1,                  0.333029, 0.9990869999999999, 0.1115571146137813,  0.1115571146137814,  14.208807503904,   127.879267535136,  0.06245467840565174, 0.06245467840565175),
    @("FAPs", "Sema6a", "Plxna4", "ECs",  3, 1,                  4.56506,            13.69518,            0.05990162668197795, 0.05990162668197797, 2, 0.6666666666666666, 2.431926, 7.295778,           0.8146397085966534, 0.8146397085966534, 11.10188810556,    99.91699295004001, 0.04879824370467203, 0.04879824370467205),
    @("FAPs", "Sema6a", "Plxna4", "FAPs", 3, 1,                  4.56506,            13.69518,            0.05990162668197795, 0.05990162668197797, 3, 1,                  0.220323, 0.6609689999999999, 0.07380317678956531, 0.07380317678956533, 1.00578771438,     9.052089429419999, 0.004420930343992562, 0.004420930343992564),
    @("FAPs", "Sema6a", "Plxna4", "sCs",  3, 1,                  4.56506,            13.69518,            0.05990162668197795, 0.05990162668197797, 3, 1,                  0.333029, 0.9990869999999999, 0.1115571146137813,  0.1115571146137814,  1.52029736674,     13.68267630066,    0.006682452633313357, 0.006682452633313359),
    @("sCs",  "Sema6a", "Plxna4", "ECs",  3, 1,                  28.97884666666667, 86.93653999999999,  0.3802535026266791, 0.3802535026266792, 2, 0.6666666666666666, 2.431926, 7.295778,           0.8146397085966534, 0.8146397085966534, 70.47441065868,    634.26969592812,   0.3097696025726546,  0.3097696025726547),
    @("sCs",  "Sema6a", "Plxna4", "FAPs", 3, 1,                  28.97884666666667, 86.93653999999999,  0.3802535026266791, 0.3802535026266792, 3, 1,                  0.220323, 0.6609689999999999, 0.07380317678956531, 0.07380317678956533, 6.384706434139999, 57.46235790725999, 0.02806391647920824, 0.02806391647920825),
    @("sCs",  "Sema6a", "Plxna4", "sCs",  3, 1,                  28.97884666666667, 86.93653999999999,  0.3802535026266791, 0.3802535026266792, 3, 1,                  0.333029, 0.9990869999999999, 0.1115571146137813,  0.1115571146137814,  9.650796326553332, 86.85716693897999, 0.04241998357481624, 0.04241998357481626)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
